# Web portal.pptx edit
#
# The underlying author edit merged a handful of two-run text fragments
# (leftovers of earlier "insert-text" / autocorrect splits) back into a
# single run on slide 3 ("Excel file") and slide 4 ("Excel Upload Form",
# "Upload Data and Metadata", "Metadata Template"). Re-assigning the full
# character range (rather than TextRange.Text, which only appends a diff
# run) collapses the paragraph back down to one run using the formatting
# of the first original run - matching how the merge looks in the
# canonical XML.

$p = $ppt.ActivePresentation

function Merge-ShapeText {
    param($shape, [string]$newText)
    $tr = $shape.TextFrame.TextRange
    $len = $tr.Length
    $full = $tr.Characters(1, $len)
    $full.Text = $newText
}

# --- Slide 3 ("MAP INTERFACE" view): "Excel " + "file" -> "Excel file" ---
$slide3 = $p.Slides.Item(3)
$excelFile = $slide3.Shapes.Item(7)
Merge-ShapeText $excelFile "Excel file"

# --- Slide 4 ("Batch Upload" view) ---
$slide4 = $p.Slides.Item(4)

# "Excel " + "Upload Form" -> "Excel Upload Form"
$excelUploadForm = $slide4.Shapes.Item(2)
Merge-ShapeText $excelUploadForm "Excel Upload Form"

# "Upload " + "Data and Metadata" -> "Upload Data and Metadata"
$uploadDataMeta = $slide4.Shapes.Item(5)
Merge-ShapeText $uploadDataMeta "Upload Data and Metadata"

# "Metad" + "ata Template" -> "Metadata Template"
$metadataTemplate = $slide4.Shapes.Item(11)
Merge-ShapeText $metadataTemplate "Metadata Template"
